$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 650.7143
$ws.Range("I20").Value = 650.7143
$ws.Range("K20").Value = 650.7143
$ws.Range("M20").Value = -420.7143
$ws.Range("H35").Value = 650.7143
$ws.Range("I35").Value = 650.7143
$ws.Range("K35").Value = 650.7143
$ws.Range("M35").Value = -271.7143
$ws.Range("H39").Value = 4132.6113
$ws.Range("I39").Value = 826.2727
$ws.Range("J39").Value = 9328.286
$ws.Range("K39").Value = 2478.8181
$ws.Range("L39").Value = 27984.858
$ws.Range("M39").Value = -2182.8181
$ws.Range("N39").Value = -28576.858
$ws.Range("H40").Value = 5638.4614
$ws.Range("I40").Value = 3830
$ws.Range("J40").Value = 11666.667
$ws.Range("K40").Value = 3830
$ws.Range("L40").Value = 11666.667
$ws.Range("M40").Value = -3655
$ws.Range("N40").Value = -12016.667
$ws.Range("H48").Value = 4450
$ws.Range("I48").Value = 5000
$ws.Range("J48").Value = 3625
$ws.Range("K48").Value = 15000
$ws.Range("L48").Value = 10875
$ws.Range("M48").Value = -14708
$ws.Range("N48").Value = -11459
$ws.Range("H51").Value = 0
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("M51").ClearContents()
$ws.Range("N51").ClearContents()
$ws.Range("H53").Value = 1692.375
$ws.Range("I53").Value = 308
$ws.Range("K53").Value = 308
$ws.Range("M53").Value = 329
$ws.Range("H56").Value = 4450
$ws.Range("I56").Value = 5000
$ws.Range("J56").Value = 3625
$ws.Range("K56").Value = 15000
$ws.Range("L56").Value = 10875
$ws.Range("M56").Value = -14466
$ws.Range("N56").Value = -11943
$ws.Range("H76").Value = 7500
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 7500
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 7500
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -8130
$ws.Range("H79").Value = 7500
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 7500
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 7500
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -9684
$ws.Range("H80").Value = 1331.8823
$ws.Range("I80").Value = 1977.3
$ws.Range("J80").Value = 409.85715
$ws.Range("K80").Value = 5931.9
$ws.Range("L80").Value = 1229.57145
$ws.Range("M80").Value = -4933.9
$ws.Range("N80").Value = -3225.57145
$ws.Range("H83").Value = 1331.8823
$ws.Range("I83").Value = 1977.3
$ws.Range("J83").Value = 409.85715
$ws.Range("K83").Value = 17795.7
$ws.Range("L83").Value = 3688.71435
$ws.Range("M83").Value = -12803.7
$ws.Range("N83").Value = -13672.71435
$ws.Range("H88").Value = 52725776
$ws.Range("J88").Value = 5957846.5
$ws.Range("L88").Value = 5957846.5
$ws.Range("N88").Value = -5958658.5
$ws.Range("H91").Value = 52725776
$ws.Range("J91").Value = 5957846.5
$ws.Range("L91").Value = 5957846.5
$ws.Range("N91").Value = -5960654.5
$ws.Range("H98").Value = 1131.4138
$ws.Range("I98").Value = 1131.4138
$ws.Range("K98").Value = 1131.4138
$ws.Range("M98").Value = 366.5862
$ws.Range("H107").Value = 810.3333
$ws.Range("I107").Value = 873.8
$ws.Range("J107").Value = 493
$ws.Range("K107").Value = 873.8
$ws.Range("L107").Value = 493
$ws.Range("M107").Value = 1046.2
$ws.Range("N107").Value = -4333
$ws.Range("H122").Value = 1131.4138
$ws.Range("I122").Value = 1131.4138
$ws.Range("K122").Value = 3394.2414
$ws.Range("M122").Value = -944.2413999999999
$ws.Range("H141").Value = 2248.7646
$ws.Range("I141").Value = 2248.7646
$ws.Range("K141").Value = 6746.293799999999
$ws.Range("M141").Value = -1566.293799999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6000
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 6000
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 6000
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -6754
$ws.Range("H61").Value = 83335976
$ws.Range("I61").Value = 100002670
$ws.Range("K61").Value = 100002670
$ws.Range("M61").Value = -100002458
$ws.Range("H74").Value = 25003362
$ws.Range("I74").Value = 29414650
$ws.Range("J74").Value = 6058.5
$ws.Range("K74").Value = 29414650
$ws.Range("L74").Value = 6058.5
$ws.Range("M74").Value = -29413776
$ws.Range("N74").Value = -7806.5
$ws.Range("H77").Value = 25003362
$ws.Range("I77").Value = 29414650
$ws.Range("J77").Value = 6058.5
$ws.Range("K77").Value = 147073250
$ws.Range("L77").Value = 30292.5
$ws.Range("M77").Value = -147068882
$ws.Range("N77").Value = -39028.5
$ws.Range("H136").Value = 83335976
$ws.Range("I136").Value = 100002670
$ws.Range("K136").Value = 300008010
$ws.Range("M136").Value = -300005460

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2351.652
$ws.Range("I105").Value = 1865.9333
$ws.Range("K105").Value = 1865.9333
$ws.Range("M105").Value = -118.9332999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 13858
$ws.Range("I39").Value = 3287.25
$ws.Range("J39").Value = 34999.5
$ws.Range("K39").Value = 3287.25
$ws.Range("L39").Value = 34999.5
$ws.Range("M39").Value = -2896.25
$ws.Range("N39").Value = -35781.5
$ws.Range("H49").Value = 13858
$ws.Range("I49").Value = 3287.25
$ws.Range("J49").Value = 34999.5
$ws.Range("K49").Value = 3287.25
$ws.Range("L49").Value = 34999.5
$ws.Range("M49").Value = -3105.25
$ws.Range("N49").Value = -35363.5
$ws.Range("H94").Value = 2161.5454
$ws.Range("J94").Value = 1940
$ws.Range("L94").Value = 1940
$ws.Range("N94").Value = -2842

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 112196.78
$ws.Range("I5").Value = 200447.4
$ws.Range("K5").Value = 601342.2
$ws.Range("M5").Value = -601230.2
$ws.Range("H46").Value = 666
$ws.Range("I46").Value = 666
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 1998
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -1907
$ws.Range("N46").ClearContents()
$ws.Range("H92").Value = 420
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 420
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 1260
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = -3756
$ws.Range("H97").Value = 267.7143
$ws.Range("I97").Value = 267.7143
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 803.1428999999999
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -307.1428999999999
$ws.Range("N97").ClearContents()
$ws.Range("H117").Value = 1524.5
$ws.Range("I117").Value = 49
$ws.Range("K117").Value = 147
$ws.Range("M117").Value = 3295
$ws.Range("H132").Value = 1789.3793
$ws.Range("I132").Value = 1013
$ws.Range("K132").Value = 9117
$ws.Range("M132").Value = -6587
$ws.Range("H135").Value = 112196.78
$ws.Range("I135").Value = 200447.4
$ws.Range("K135").Value = 1804026.6
$ws.Range("M135").Value = -1801491.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6110.778
$ws.Range("I70").Value = 6466
$ws.Range("K70").Value = 6466
$ws.Range("M70").Value = -6196
$ws.Range("H73").Value = 6110.778
$ws.Range("I73").Value = 6466
$ws.Range("K73").Value = 6466
$ws.Range("M73").Value = -5530
$ws.Range("H80").Value = 3295.4
$ws.Range("I80").Value = 3463.125
$ws.Range("K80").Value = 3463.125
$ws.Range("M80").Value = -2465.125
$ws.Range("H83").Value = 3295.4
$ws.Range("I83").Value = 3463.125
$ws.Range("K83").Value = 17315.625
$ws.Range("M83").Value = -12323.625
$ws.Range("H102").Value = 4999.75
$ws.Range("I102").Value = 4999.75
$ws.Range("K102").Value = 4999.75
$ws.Range("M102").Value = -3377.75
$ws.Range("H107").Value = 2775.5
$ws.Range("J107").Value = 4951.5
$ws.Range("L107").Value = 4951.5
$ws.Range("N107").Value = -8791.5
$ws.Range("H113").Value = 29230
$ws.Range("I113").Value = 34066.527
$ws.Range("J113").Value = 4356.4287
$ws.Range("K113").Value = 34066.527
$ws.Range("L113").Value = 4356.4287
$ws.Range("M113").Value = -31896.527
$ws.Range("N113").Value = -8696.4287
$ws.Range("H124").Value = 103198
$ws.Range("J124").Value = 103198
$ws.Range("L124").Value = 103198
$ws.Range("N124").Value = -113018
$ws.Range("H139").Value = 129781
$ws.Range("J139").Value = 129781
$ws.Range("L139").Value = 129781
$ws.Range("N139").Value = -140061

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 375.66666
$ws.Range("I55").Value = 304.8
$ws.Range("K55").Value = 304.8
$ws.Range("M55").Value = -131.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 10871296
$ws.Range("I136").Value = 11365414
$ws.Range("K136").Value = 34096242
$ws.Range("M136").Value = -34093692

